# Extend the corr_lambdas matrix from 33x33 (A1:AH34) to 35x35 (A1:AJ36)
# and update the diagonal "lambda" values per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend header row (row 1) with two new index columns: AI1=33, AJ1=34 ---
$ws.Range("AH1").Copy() | Out-Null
$ws.Range("AI1:AJ1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 35).Value = 33
$ws.Cells.Item(1, 36).Value = 34

# --- 2. Extend header column (col A) with two new index rows: A35=33, A36=34 ---
$ws.Range("A34").Copy() | Out-Null
$ws.Range("A35:A36").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(36, 1).Value = 34

# --- 3. Update the diagonal lambda values for existing rows 2..34 ---
$ws.Cells.Item(2, 2).Value = 14.48433289256554
$ws.Cells.Item(3, 3).Value = 10.0736367654875
$ws.Cells.Item(4, 4).Value = 3.17659814331083
$ws.Cells.Item(5, 5).Value = 2.214924190957984
$ws.Cells.Item(6, 6).Value = 1.241904923035922
$ws.Cells.Item(7, 7).Value = 0.8215235589559928
$ws.Cells.Item(8, 8).Value = 0.885270365287341
$ws.Cells.Item(9, 9).Value = 0.09061962104880264
$ws.Cells.Item(10, 10).Value = 0.5435170356649683
$ws.Cells.Item(11, 11).Value = 0.507329304277006
$ws.Cells.Item(12, 12).Value = 0.2391122758873609
$ws.Cells.Item(13, 13).Value = 0.3281475414741265
$ws.Cells.Item(14, 14).Value = 0.3930833820466329
$ws.Cells.Item(15, 15).Value = [double]"-7.190621613715121e-16"
$ws.Cells.Item(16, 16).Value = [double]"-7.190621613715121e-16"
$ws.Cells.Item(17, 17).Value = [double]"5.698258448364944e-16"
$ws.Cells.Item(18, 18).Value = [double]"5.698258448364944e-16"
$ws.Cells.Item(19, 19).Value = [double]"-6.570937268611783e-17"
$ws.Cells.Item(20, 20).Value = [double]"-6.570937268611783e-17"
$ws.Cells.Item(21, 21).Value = [double]"3.648468800192603e-16"
$ws.Cells.Item(22, 22).Value = [double]"3.648468800192603e-16"
$ws.Cells.Item(23, 23).Value = [double]"3.175963566890125e-16"
$ws.Cells.Item(24, 24).Value = [double]"3.175963566890125e-16"
$ws.Cells.Item(25, 25).Value = [double]"-4.62986842669504e-16"
$ws.Cells.Item(26, 26).Value = [double]"-4.62986842669504e-16"
$ws.Cells.Item(27, 27).Value = [double]"1.077444818183858e-16"
$ws.Cells.Item(28, 28).Value = [double]"1.077444818183858e-16"
$ws.Cells.Item(29, 29).Value = [double]"6.53128183089937e-17"
$ws.Cells.Item(30, 30).Value = [double]"6.53128183089937e-17"
$ws.Cells.Item(31, 31).Value = [double]"-2.219548042909112e-16"
$ws.Cells.Item(32, 32).Value = [double]"-2.219548042909112e-16"
$ws.Cells.Item(33, 33).Value = [double]"-6.484476954447993e-17"
$ws.Cells.Item(34, 34).Value = [double]"-6.484476954447993e-17"

# --- 4. Fill the two new columns (AI=35, AJ=36) with 0 for existing rows 2..34 ---
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 35).Value = 0
    $ws.Cells.Item($r, 36).Value = 0
}

# --- 5. Fill the two new rows (35, 36) with 0 across columns B..AJ (2..36) ---
for ($c = 2; $c -le 36; $c++) {
    $ws.Cells.Item(35, $c).Value = 0
    $ws.Cells.Item(36, $c).Value = 0
}

# --- 6. Set the diagonal values for the two new rows: AI35 and AJ36 ---
$ws.Cells.Item(35, 35).Value = [double]"-2.954486370814103e-16"
$ws.Cells.Item(36, 36).Value = [double]"-2.43362683501509e-16"
